# Update stimulus presentation time-logging identifiers in filenames and
# sheet names (the numeric timestamps embedded in the names/files).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555013240163"
$ws1.Range("B2").Value = "go_stims-16512555012960193.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255501307019.csv"
$ws1.Range("B4").Value = "go_stims-16512555013090208.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555013230505.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512555033068933"
$ws2.Range("B2").Value = "ZB-match_1-1651255501525297.csv"
$ws2.Range("B3").Value = "OB-16512555017833283.csv"
$ws2.Range("B4").Value = "TB-16512555022123375.csv"
$ws2.Range("B5").Value = "OB-16512555016763012.csv"
$ws2.Range("B6").Value = "TB-16512555032944932.csv"
$ws2.Range("B7").Value = "TB-16512555031964931.csv"
$ws2.Range("B8").Value = "ZB-match_8-16512555014763405.csv"
$ws2.Range("B9").Value = "ZB-match_3-16512555014199767.csv"
$ws2.Range("B10").Value = "OB-16512555016512978.csv"

# --- Sheet 3: RS_TO (only the sheet name changes) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555033138947"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651255503372893"
$ws4.Range("B2").Value = "MM_stims-16512555033388927.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555033158956.csv"
$ws4.Range("B4").Value = "MM_stims-16512555033548956.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255503339894.csv"
$ws4.Range("B6").Value = "MM_stims-16512555033708935.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555033558972.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16512555034508934"
$ws5.Range("B2").Value = "SAT_stims-16512555034038935.csv"
$ws5.Range("B3").Value = "vSAT_stims-165125550341995.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555033778942.csv"
$ws5.Range("B5").Value = "vSAT_stims-16512555034348927.csv"
